$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new time-tracking entries for row 11 (copy formatting down from row 10)
$ws.Range("B10:F10").Copy() | Out-Null
$ws.Range("B11:F11").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B11").Value = 44824
$ws.Range("C11").Value = 0.48958333333333331
$ws.Range("D11").Value = 0.53125
$ws.Range("E11").Formula = "=D11-C11"

# Row 12
$ws.Range("B12").Value = 44824
$ws.Range("C12").Value = 0.5625
$ws.Range("D12").Value = 0.63541666666666663
$ws.Range("E12").Formula = "=D12-C12"

# Update the running total formula in F7 to include the newly-used rows through 23
$ws.Range("F7").Formula = "=SUM(E7:E23)"

# Note progress in the comment column of row 12
$ws.Range("G12").Value = "Finished upto lesson 27"

# Update the active selection to match the author's final cursor position
$ws.Range("G13").Select()

$wb.Application.CalculateFull()
